$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.265.92"
$ws.Range("E2").Value = "  -0.25%  "

# Row 3
$ws.Range("D3").Value = "1.929.38"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'249.09"
$ws.Range("E5").Value = "  -0.30%  "

# Row 6
$ws.Range("D6").Value = "'0.7170"
$ws.Range("E6").Value = "  -0.94%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.3203"
$ws.Range("E8").Value = "  -4.36%  "

# Row 9
$ws.Range("D9").Value = "'27.67"
$ws.Range("E9").Value = "  -2.83%  "

# Row 10
$ws.Range("D10").Value = "'0.07103"
$ws.Range("E10").Value = "  -3.44%  "

# Row 11
$ws.Range("D11").Value = "'0.7914"
$ws.Range("E11").Value = "  -2.94%  "

# Row 12
$ws.Range("D12").Value = "'0.07992"
$ws.Range("E12").Value = "  -1.74%  "

# Row 13
$ws.Range("D13").Value = "1.932.32"
$ws.Range("E13").Value = "  -0.31%  "

# Row 14
$ws.Range("D14").Value = "'5.395"
$ws.Range("E14").Value = "  -2.39%  "

# Row 15
$ws.Range("D15").Value = "'94.85"
$ws.Range("E15").Value = "  -0.44%  "

# Row 16
$ws.Range("D16").Value = "'14.67"
$ws.Range("E16").Value = "  -1.49%  "

# Row 17
$ws.Range("D17").Value = "30.263.87"
$ws.Range("E17").Value = "  -0.24%  "

# Row 18
$ws.Range("D18").Value = "'256.75"
$ws.Range("E18").Value = "  +0.81%  "

# Row 19
$ws.Range("D19").Value = "'0.000008039"
$ws.Range("E19").Value = "  -3.78%  "

# Row 20
$ws.Range("D20").Value = "'5.771"
$ws.Range("E20").Value = "  -1.66%  "

# Row 21
$ws.Range("D21").Value = "2.178.95"
$ws.Range("E21").Value = "  -0.37%  "

# Row 22
$ws.Range("D22").Value = "'0.9993"
$ws.Range("E22").Value = "  -0.06%  "

# Row 23
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").Value = "'6.824"
$ws.Range("E24").Value = "  -1.88%  "

# Row 25
$ws.Range("D25").Value = "'9.546"
$ws.Range("E25").Value = "  -2.91%  "

# Row 26
$ws.Range("D26").Value = "'164.93"
$ws.Range("E26").Value = "  +2.75%  "

# Row 27
$ws.Range("E27").Value = "  -2.78%  "

# Row 28
$ws.Range("D28").Value = "'2.270"
$ws.Range("E28").Value = "  -6.47%  "

# Row 29
$ws.Range("D29").Value = "'0.1269"
$ws.Range("E29").Value = "  -4.65%  "

# Row 30
$ws.Range("E30").Value = "  +0.87%  "

# Row 31
$ws.Range("D31").Value = "'1.526"
$ws.Range("E31").Value = "  -2.29%  "

# Row 32
$ws.Range("E32").Value = "  -1.21%  "

# Row 33
$ws.Range("D33").Value = "'4.129"
$ws.Range("E33").Value = "  -2.82%  "

# Row 34
$ws.Range("D34").Value = "'0.05141"
$ws.Range("E34").Value = "  -1.64%  "

# Row 35
$ws.Range("D35").Value = "'1.268"
$ws.Range("E35").Value = "  +0.58%  "

# Row 36
$ws.Range("E36").Value = "  -0.89%  "

# Row 37
$ws.Range("D37").Value = "'2.763"
$ws.Range("E37").Value = "  +0.90%  "

# Row 38
$ws.Range("E38").Value = "  -1.97%  "

# Row 39
$ws.Range("E39").Value = "  -1.69%  "

# Row 40
$ws.Range("D40").Value = "'79.34"
$ws.Range("E40").Value = "  -0.60%  "

# Row 41
$ws.Range("D41").Value = "'6.358"
$ws.Range("E41").Value = "  -4.85%  "

# Row 42
$ws.Range("E42").Value = "  -1.00%  "

# Row 43
$ws.Range("D43").Value = "'1.994"
$ws.Range("E43").Value = "  -1.74%  "

# Row 44
$ws.Range("D44").Value = "'0.8481"
$ws.Range("E44").Value = "  +0.63%  "

# Row 45
$ws.Range("D45").Value = "'0.9993"
$ws.Range("E45").Value = "  -0.13%  "

# Row 46
$ws.Range("D46").Value = "'100.55"
$ws.Range("E46").Value = "  -2.27%  "

# Row 47
$ws.Range("D47").Value = "'9.759"
$ws.Range("E47").Value = "  -0.56%  "

# Row 48
$ws.Range("D48").Value = "'7.433"
$ws.Range("E48").Value = "  +0.33%  "

# Row 49
$ws.Range("D49").Value = "'36.68"
$ws.Range("E49").Value = "  -0.56%  "

# Row 50
$ws.Range("D50").Value = "'957.01"
$ws.Range("E50").Value = "  +11.08%  "

# Row 51
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "'0.4204"
$ws.Range("E51").Value = "  +1.18%  "
